$d = $word.ActiveDocument

# [4] "add hero credit" — strike through the whole paragraph, including the
# paragraph mark itself (so both the run and the pPr's rPr pick up <w:strike/>).
$p4 = $d.Paragraphs(2)
$p4.Range.Font.StrikeThrough = 1

# [8] "toggle the visibility of the favorites panel" — only the first run of
# this paragraph is struck through; the trailing " with transition" run and
# the paragraph mark are left alone, so select just that run's text.
$r8 = $d.Content
$found8 = $r8.Find.Execute("[8] toggle the visibility of the favorites panel", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r8.Font.StrikeThrough = 1

# [11] "Play Details Missing everything" — strike through the whole
# paragraph (all runs + the paragraph mark's rPr).
$p11 = $d.Paragraphs(14)
$p11.Range.Font.StrikeThrough = 1
